$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 570-571, shifting the existing
# data (rows 570-676) down to rows 572-678.
$ws.Range("A570:R571").Insert()

# Fill in the two newly inserted rows with the new weekly readings.
# Row 570 - "Primera" quality
$ws.Cells.Item(570, 1).Value = 3
$ws.Cells.Item(570, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(570, 3).Value = "Coquimbo"
$ws.Cells.Item(570, 4).Value = 44694
$ws.Cells.Item(570, 5).Value = 5
$ws.Cells.Item(570, 6).Value = 100114014
$ws.Cells.Item(570, 7).Value = "Betarraga"
$ws.Cells.Item(570, 8).Value = "Sin especificar"
$ws.Cells.Item(570, 9).Value = "Primera"
$ws.Cells.Item(570, 10).Value = 2500
$ws.Cells.Item(570, 11).Value = 650
$ws.Cells.Item(570, 12).Value = 700
$ws.Cells.Item(570, 13).Value = 676
$ws.Cells.Item(570, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(570, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(570, 16).Value = 169
$ws.Cells.Item(570, 17).Value = 4
$ws.Cells.Item(570, 18).Value = "Hortaliza"

# Row 571 - "Segunda" quality
$ws.Cells.Item(571, 1).Value = 3
$ws.Cells.Item(571, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(571, 3).Value = "Coquimbo"
$ws.Cells.Item(571, 4).Value = 44694
$ws.Cells.Item(571, 5).Value = 5
$ws.Cells.Item(571, 6).Value = 100114014
$ws.Cells.Item(571, 7).Value = "Betarraga"
$ws.Cells.Item(571, 8).Value = "Sin especificar"
$ws.Cells.Item(571, 9).Value = "Segunda"
$ws.Cells.Item(571, 10).Value = 1400
$ws.Cells.Item(571, 11).Value = 450
$ws.Cells.Item(571, 12).Value = 450
$ws.Cells.Item(571, 13).Value = 450
$ws.Cells.Item(571, 14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(571, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(571, 16).Value = 112
$ws.Cells.Item(571, 17).Value = 4
$ws.Cells.Item(571, 18).Value = "Hortaliza"
